# Update "想去人数" (F column) figures on the 展览, 演出, and 全部类型 sheets
# to reflect the latest output snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (sheetId 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4746
$ws1.Range("F4").Value = 3207
$ws1.Range("F6").Value = 562
$ws1.Range("F7").Value = 567
$ws1.Range("F8").Value = 422
$ws1.Range("F9").Value = 146
$ws1.Range("F10").Value = 1814
$ws1.Range("F11").Value = 1406
$ws1.Range("F13").Value = 1670
$ws1.Range("F14").Value = 27
$ws1.Range("F16").Value = 633
$ws1.Range("F19").Value = 544
$ws1.Range("F21").Value = 63
$ws1.Range("F22").Value = 125
$ws1.Range("F23").Value = 14
$ws1.Range("F24").Value = 119
$ws1.Range("F25").Value = 52
$ws1.Range("F26").Value = 93
$ws1.Range("F27").Value = 4182
$ws1.Range("F31").Value = 1959
$ws1.Range("F33").Value = 1954

# --- Sheet: 演出 (sheetId 2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 62

# --- Sheet: 全部类型 (sheetId 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4746
$ws4.Range("F4").Value = 3207
$ws4.Range("F6").Value = 562
$ws4.Range("F7").Value = 567
$ws4.Range("F9").Value = 422
$ws4.Range("F10").Value = 146
$ws4.Range("F11").Value = 1814
$ws4.Range("F12").Value = 1406
$ws4.Range("F14").Value = 1670
$ws4.Range("F15").Value = 27
$ws4.Range("F17").Value = 633
$ws4.Range("F20").Value = 544
$ws4.Range("F22").Value = 63
$ws4.Range("F23").Value = 125
$ws4.Range("F24").Value = 14
$ws4.Range("F25").Value = 119
$ws4.Range("F26").Value = 52
$ws4.Range("F27").Value = 93
$ws4.Range("F28").Value = 4182
$ws4.Range("F29").Value = 62
$ws4.Range("F34").Value = 1959
$ws4.Range("F36").Value = 1954
